$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "2026-01-12 14:57:58"
